# Scheduled market-data refresh: Universalis current-price pulls changed, so the
# derived currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H,I,J,K,L,M,N) are recomputed per-row on every affected Leve sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1400
$ws.Range("I2").Value = 1400
$ws.Range("K2").Value = 1400
$ws.Range("M2").Value = -1287
$ws.Range("H75").Value = 252445
$ws.Range("J75").Value = 252445
$ws.Range("L75").Value = 252445
$ws.Range("N75").Value = -254317
$ws.Range("H78").Value = 252445
$ws.Range("J78").Value = 252445
$ws.Range("L78").Value = 757335
$ws.Range("N78").Value = -766695
$ws.Range("H88").Value = 8549232
$ws.Range("I88").Value = 1566.6666
$ws.Range("J88").Value = 11113532
$ws.Range("K88").Value = 1566.6666
$ws.Range("L88").Value = 11113532
$ws.Range("M88").Value = -1160.6666
$ws.Range("N88").Value = -11114344
$ws.Range("H91").Value = 8549232
$ws.Range("I91").Value = 1566.6666
$ws.Range("J91").Value = 11113532
$ws.Range("K91").Value = 1566.6666
$ws.Range("L91").Value = 11113532
$ws.Range("M91").Value = -162.6666
$ws.Range("N91").Value = -11116340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2580.9468
$ws.Range("I32").Value = 1717.85
$ws.Range("K32").Value = 1717.85
$ws.Range("M32").Value = -1430.85
$ws.Range("H61").Value = 2897.077
$ws.Range("I61").Value = 1933.4667
$ws.Range("K61").Value = 1933.4667
$ws.Range("M61").Value = -1721.4667
$ws.Range("H74").Value = 1244.6578
$ws.Range("I74").Value = 1297.32
$ws.Range("J74").Value = 1143.3846
$ws.Range("K74").Value = 1297.32
$ws.Range("L74").Value = 1143.3846
$ws.Range("M74").Value = -423.3199999999999
$ws.Range("N74").Value = -2891.3846
$ws.Range("H77").Value = 1244.6578
$ws.Range("I77").Value = 1297.32
$ws.Range("J77").Value = 1143.3846
$ws.Range("K77").Value = 6486.599999999999
$ws.Range("L77").Value = 5716.923000000001
$ws.Range("M77").Value = -2118.599999999999
$ws.Range("N77").Value = -14452.923
$ws.Range("H109").Value = 193459
$ws.Range("J109").Value = 193459
$ws.Range("L109").Value = 193459
$ws.Range("N109").Value = -196233
$ws.Range("H122").Value = 13802.875
$ws.Range("I122").Value = 13802.875
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 41408.625
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -38958.625
$ws.Range("N122").Value = $null
$ws.Range("H123").Value = 39714.5
$ws.Range("J123").Value = 39714.5
$ws.Range("L123").Value = 39714.5
$ws.Range("N123").Value = -49514.5
$ws.Range("H136").Value = 2897.077
$ws.Range("I136").Value = 1933.4667
$ws.Range("K136").Value = 5800.4001
$ws.Range("M136").Value = -3250.4001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
$ws.Range("H134").Value = 2211.1345
$ws.Range("I134").Value = 1389.4651
$ws.Range("J134").Value = 6136.8887
$ws.Range("K134").Value = 4168.3953
$ws.Range("L134").Value = 18410.6661
$ws.Range("M134").Value = -1633.3953
$ws.Range("N134").Value = -23480.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 49537
$ws.Range("J68").Value = 49537
$ws.Range("L68").Value = 49537
$ws.Range("N68").Value = -51035
$ws.Range("H71").Value = 49537
$ws.Range("J71").Value = 49537
$ws.Range("L71").Value = 148611
$ws.Range("N71").Value = -156099

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1882.2858
$ws.Range("I34").Value = 488
$ws.Range("J34").Value = 2440
$ws.Range("K34").Value = 1464
$ws.Range("L34").Value = 7320
$ws.Range("M34").Value = -1380
$ws.Range("N34").Value = -7488
$ws.Range("H39").Value = 8931.6
$ws.Range("J39").Value = 8931.6
$ws.Range("L39").Value = 26794.8
$ws.Range("N39").Value = -27382.8
$ws.Range("H107").Value = 983.678
$ws.Range("I107").Value = 332.5
$ws.Range("J107").Value = 1225.9767
$ws.Range("K107").Value = 997.5
$ws.Range("L107").Value = 3677.9301
$ws.Range("M107").Value = 922.5
$ws.Range("N107").Value = -7517.9301

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 22810.385
$ws.Range("J93").Value = 22810.385
$ws.Range("L93").Value = 22810.385
$ws.Range("N93").Value = -26554.385
$ws.Range("H102").Value = 6439.1665
$ws.Range("I102").Value = 1878.3334
$ws.Range("K102").Value = 1878.3334
$ws.Range("M102").Value = -256.3334
$ws.Range("H132").Value = 4135.24
$ws.Range("I132").Value = 3807.889
$ws.Range("J132").Value = 4977
$ws.Range("K132").Value = 11423.667
$ws.Range("L132").Value = 14931
$ws.Range("M132").Value = -8893.667000000001
$ws.Range("N132").Value = -19991

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1052.8572
$ws.Range("I22").Value = 1078
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 1078
$ws.Range("L22").Value = 990
$ws.Range("M22").Value = -783
$ws.Range("N22").Value = -1580
$ws.Range("H27").Value = 1052.8572
$ws.Range("I27").Value = 1078
$ws.Range("J27").Value = 990
$ws.Range("K27").Value = 1078
$ws.Range("L27").Value = 990
$ws.Range("M27").Value = -971
$ws.Range("N27").Value = -1204
$ws.Range("H62").Value = 48950
$ws.Range("J62").Value = 48950
$ws.Range("L62").Value = 48950
$ws.Range("N62").Value = -50198
$ws.Range("H65").Value = 48950
$ws.Range("J65").Value = 48950
$ws.Range("L65").Value = 146850
$ws.Range("N65").Value = -153090
$ws.Range("H125").Value = 46800
$ws.Range("J125").Value = 46800
$ws.Range("L125").Value = 46800
$ws.Range("N125").Value = -56640
$ws.Range("H136").Value = 5024.8667
$ws.Range("I136").Value = 1855.7778
$ws.Range("J136").Value = 9778.5
$ws.Range("K136").Value = 5567.3334
$ws.Range("L136").Value = 29335.5
$ws.Range("M136").Value = -3017.3334
$ws.Range("N136").Value = -34435.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15153468
$ws.Range("I132").Value = 17858556
$ws.Range("K132").Value = 53575668
$ws.Range("M132").Value = -53573138
$ws.Range("H136").Value = 13374630
$ws.Range("I136").Value = 15921252
$ws.Range("K136").Value = 47763756
$ws.Range("M136").Value = -47761206

